# Update the regression-equation labels for the "maximum" series (J13) and
# the average/other series (M13) to reflect the refreshed fit once row 9-12
# data changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mapping")

$ws.Range("J13").Value = "y = 1.6304x - 26.059"
$ws.Range("M13").Value = "y = 1.7566x - 34.321"

# Insert an XY scatter chart plotting the "maximum" bucket fit (columns L/M,
# rows 9-12) with a linear trendline whose equation is displayed on the
# chart, mirroring the chart that was added next to the data table.
$chartObj = $ws.Shapes.AddChart2(-1, -4169)
$chart = $chartObj.Chart
$chart.ChartType = -4169

$series = $chart.SeriesCollection().NewSeries()
$series.XValues = "=mapping!`$L`$9:`$L`$12"
$series.Values = "=mapping!`$M`$9:`$M`$12"
$series.MarkerStyle = 8
$series.MarkerSize = 5

$trend = $series.Trendlines().Add()
$trend.Type = -4132
$trend.DisplayEquation = $true
$trend.DisplayRSquared = $false

$chart.HasLegend = $false
$chart.HasTitle = $false

# Position the chart roughly over G25:O40, matching where it was dropped on
# the sheet.
$chartObj.Left = $ws.Range("G25").Left
$chartObj.Top = $ws.Range("G25").Top
$chartObj.Width = $ws.Range("G25:N39").Width
$chartObj.Height = $ws.Range("G25:N39").Height

# Restore the view so the window is scrolled to show column F onward with
# O19 as the active cell, matching the saved sheet view.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("O19").Select()

$wb.Save()
